$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "26-03-2025"
$ws.Range("B7").Value = "Rajasthan Royals vs Kolkata Knight Riders"
$ws.Range("C7").Value = "Kolkata Knight Riders"
$ws.Range("D7").Value = "Kolkata Knight Riders"
